$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.456.97"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "3.692.67"
$ws.Range("E3").Value = "  -3.17%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'687.42"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "'160.50"
$ws.Range("E6").Value = "  -5.97%  "
$ws.Range("D7").Value = "3.691.18"
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -5.92%  "
$ws.Range("E10").Value = "  -8.90%  "
$ws.Range("D11").Value = "'7.21"
$ws.Range("E11").Value = "  -4.06%  "
$ws.Range("E12").Value = "  -10.02%  "
$ws.Range("E13").Value = "  -7.13%  "
$ws.Range("D14").Value = "4.317.72"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "'32.60"
$ws.Range("E15").Value = "  -10.45%  "
$ws.Range("D16").Value = "3.702.99"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").Value = "69.457.76"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'15.98"
$ws.Range("E19").Value = "  -9.52%  "
$ws.Range("D20").Value = "'6.46"
$ws.Range("E20").Value = "  -10.84%  "
$ws.Range("D21").Value = "'472.90"
$ws.Range("E21").Value = "  -7.94%  "
$ws.Range("D22").Value = "'9.95"
$ws.Range("E22").Value = "  -5.40%  "
$ws.Range("D23").Value = "'0.648"
$ws.Range("E23").Value = "  -9.58%  "
$ws.Range("D24").Value = "'79.78"
$ws.Range("E24").Value = "  -4.70%  "
$ws.Range("D25").Value = "3.840.03"
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'0.0000125"
$ws.Range("E27").Value = "  -11.80%  "
$ws.Range("E28").Value = "  -13.55%  "
$ws.Range("D29").Value = "'9.23"
$ws.Range("E29").Value = "  -10.50%  "
$ws.Range("D30").Value = "'2.71"
$ws.Range("E30").Value = "  -9.80%  "
$ws.Range("D31").Value = "'1.76"
$ws.Range("E31").Value = "  -12.42%  "
$ws.Range("D32").Value = "'6.68"
$ws.Range("E32").Value = "  -8.84%  "
$ws.Range("E33").Value = "  -11.25%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'26.76"
$ws.Range("E35").Value = "  -8.46%  "
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("D37").Value = "'8.21"
$ws.Range("E37").Value = "  -12.22%  "
$ws.Range("D38").Value = "'6.16"
$ws.Range("E38").Value = "  -7.78%  "
$ws.Range("D39").Value = "'2.29"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'0.0908"
$ws.Range("E41").Value = "  -10.17%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'167.82"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").Value = "'0.943"
$ws.Range("E44").Value = "  -6.88%  "
$ws.Range("D45").Value = "'47.89"
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("D46").Value = "'2.74"
$ws.Range("E46").Value = "  -15.46%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.31"
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("B48").Value = "SuiNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").Value = "'1.12"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("E49").Value = "  -9.19%  "
$ws.Range("D50").Value = "'28.47"
$ws.Range("E50").Value = "  -6.05%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'374.77"
$ws.Range("E51").Value = "  -13.54%  "
